$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before old row 12 (so old row 11 stays, and rows 11-16 become the new block)
$ws.Rows("12:16").Insert()

# --- Text values are entered in a specific order so shared-string indices line up ---

# A11
$ws.Range("A11").Value = "Test the WebServices Eric created"

# C11
$ws.Range("C11").Value = "Dave - worked great, but thought we should add methods to return a String.  I also created the repo to test the String results and used the built-in jersey server. Put these in EmailValidationRestTest."
$ws.Range("C11").WrapText = $true

# C13 hyperlink (email) - creates the Hyperlink style
$ws.Range("C13").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:dbsullivan@madisoncollege", "", "", "dbsullivan@madisoncollege ")

# C12
$ws.Range("C12").Value = "Test result Differences noted:"
$ws.Range("C12").WrapText = $true

# D12 / E12 (bold + wrap) - creates the bold-wrap style
$ws.Range("D12").Value = "isEmailValid (InternetAddress)"
$ws.Range("D12").Font.Bold = $true
$ws.Range("D12").WrapText = $true
$ws.Range("E12").Value = "isEmailValid2 (RegEx)"
$ws.Range("E12").Font.Bold = $true
$ws.Range("E12").WrapText = $true

# A13
$ws.Range("A13").Value = "When I setup tests, I noticed user might want to combine tests to avoid conflicting results. Regex might be more accurate when extension is omitted."
$ws.Range("A13").WrapText = $true

# A15
$ws.Range("A15").Value = "I plan to put this into my TennisApp as an edit on Player entry"

# --- Non-string values / formats ---

$ws.Range("B11").Value = 42303
$ws.Range("B12").Value = 42303
$ws.Range("B13").Value = 42303
$ws.Range("B15").Value = 42310

$ws.Range("D13").Value = $true
$ws.Range("E13").Value = $false

$ws.Range("C14").WrapText = $true
$ws.Range("C15").WrapText = $true

# Row heights
$ws.Rows("11").RowHeight = 60
$ws.Rows("12").RowHeight = 30
$ws.Rows("13").RowHeight = 45

# Column widths (engine quantizes ColumnWidth to whole pixels, so these are the
# closest achievable inputs to the target XML widths of 19.85546875 / 19.5703125)
$ws.Columns("D").ColumnWidth = 19
$ws.Columns("E").ColumnWidth = 18.6

# Update selection to match target
$ws.Range("C11").Select()
